$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $wb.Windows.Item(1)

# New rows appended to the variable_definitions table (rows 134-142).
# Columns A/B/C/E/F follow the same constant pattern as the preceding
# "Emissions" rows (NA, Emissions, NA, ..., Mt CO2/yr, created for LEEP
# report data). Column D (the Variable name) must be written in the
# specific order below so that new shared-string entries are appended
# in the same order as the source workbook.
$ws.Cells.Item(134, 4).Value = "Emissions|CO2|Energy|Demand|Industry|Electricity"
$ws.Cells.Item(136, 4).Value = "Emissions|CO2|Energy|Demand|Buildings|Electricity"
$ws.Cells.Item(135, 4).Value = "Emissions|CO2|Energy|Demand|Transportation|Electricity"
$ws.Cells.Item(137, 4).Value = "Emissions|CO2|Energy|Demand|Industry|TotalwElec"
$ws.Cells.Item(138, 4).Value = "Emissions|CO2|Energy|Demand|Transportation|TotalwElec"
$ws.Cells.Item(139, 4).Value = "Emissions|CO2|Energy|Demand|Buildings|TotalwElec"
$ws.Cells.Item(140, 4).Value = "Emissions|CO2|Energy|Demand|Industry|TotalDI"
$ws.Cells.Item(141, 4).Value = "Emissions|CO2|Energy|Demand|Transportation|TotalDI"
$ws.Cells.Item(142, 4).Value = "Emissions|CO2|Energy|Demand|Buildings|TotalDI"

for ($row = 134; $row -le 142; $row++) {
    $ws.Cells.Item($row, 1).Value = "NA"
    $ws.Cells.Item($row, 2).Value = "Emissions"
    $ws.Cells.Item($row, 3).Value = "NA"
    $ws.Cells.Item($row, 5).Value = "Mt CO2/yr"
    $ws.Cells.Item($row, 6).Value = "created for LEEP report data"
}

# Update the window/selection state to match: scrolled to the new bottom
# rows, with D143 (the next empty row in column D) selected.
$ws.Activate()
$ws.Range("D143").Select()
$win.ScrollRow = 128
$win.ScrollColumn = 1
